$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet
$ws.Name = "Physics Grades"

# Update row 2: B2 becomes a number, C2 becomes "F", D2 becomes "Fail"
$ws.Range("B2").Value = 123
$ws.Range("C2").Value = "F"
$ws.Range("D2").Value = "Fail"

# Delete row 3 entirely
$ws.Rows("3:3").Delete()
